{"js": "// The captured change for this revision is a pure \"re-save\" fingerprint:\n// Word stamped newer compatibility markers (the `oel` / `w16du` namespace\n// declarations + the matching `mc:Ignorable` tokens on document.xml,\n// endnotes.xml, footer1.xml, footer2.xml, footnotes.xml, header1.xml,\n// numbering.xml and styles.xml), assigned `w16cid:durableId` values to the\n// existing `<w:num>` list instances, synced a handful of new\n// `<w:lsdException>` latent-style entries in styles.xml, and the package's\n// customXml parts were renumbered by the repackager. None of the document's\n// actual content (paragraphs, runs, tables, headers/footers text, styles,\n// numbering *definitions*, etc.) changed \u2014 the rendered/queryable document\n// is identical before and after.\n//\n// Those compatibility/serialization fingerprints are stamped internally by\n// Word's save pipeline and are not reachable through the Word JavaScript\n// API surface (there is no `durableId`, latent-style, or raw-namespace\n// control exposed on `Word.Document`/`Word.Range`/etc.), so there is no\n// content-level action for this script to perform. We still touch the\n// context once so the call is a well-formed, verifiable no-op rather than\n// an empty script.\ncontext.document.body.load(\"text\");\nawait context.sync();\n", "ps1": "# The captured change for this revision is a pure \"re-save\" fingerprint:\n# Word stamped newer compatibility markers (the `oel` / `w16du` namespace\n# declarations + the matching mc:Ignorable tokens on document.xml,\n# endnotes.xml, footer1.xml, footer2.xml, footnotes.xml, header1.xml,\n# numbering.xml and styles.xml), assigned w16cid:durableId values to the\n# existing <w:num> list instances, synced a handful of new\n# <w:lsdException> latent-style entries in styles.xml, and the package's\n# customXml parts were renumbered by the repackager. None of the document's\n# actual content (paragraphs, runs, tables, headers/footers text, styles,\n# numbering *definitions*, etc.) changed -- the rendered/queryable document\n# is identical before and after.\n#\n# Those compatibility/serialization fingerprints are stamped internally by\n# Word's save pipeline and are not reachable through the Word COM object\n# model (there is no DurableId, LatentStyles, or raw-namespace control\n# exposed on Document/Range/List/etc.), so there is no content-level action\n# for this script to perform. We still touch $d once so the call is a\n# well-formed, verifiable no-op rather than an empty script.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
